$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing shared strings in-place (keeps their shared-string indices)
$ws.Range("B2").Value = "RF 100, null=-1, no class balancing"
$ws.Range("B3").Value = "RF 100, null=median, no class balancing"
$ws.Range("B4").Value = "RF 100, null=mean, no class balancing"
$ws.Range("B5").Value = "RF 100, null managed with RF, no class balancing"

# Row 6
$ws.Range("A6").Value = 4.1
$ws.Range("B6").Value = "RF 100, null=-1, balancing = 1"
$ws.Range("C6").Value = 0.851611969901
$ws.Range("C6").NumberFormat = "0.000000"
$ws.Range("C6").WrapText = $true
$ws.Range("D6").Value = 0.00372126517529
$ws.Range("D6").NumberFormat = "0.000000"
$ws.Range("E6").Value = 0.863359
$ws.Range("E6").WrapText = $true

# Row 7
$ws.Range("A7").Value = 4.2
$ws.Range("B7").Value = "RF 100, null=-1, balancing = 2"
$ws.Range("C7").Value = 0.851354124315
$ws.Range("C7").NumberFormat = "0.000000"
$ws.Range("C7").WrapText = $true
$ws.Range("D7").Value = 0.00380425264735
$ws.Range("D7").NumberFormat = "0.000000"

# Row 8
$ws.Range("A8").Value = 4.3
$ws.Range("B8").Value = "RF 100, null=-1, balancing = 5"
$ws.Range("C8").Value = 0.849760872926
$ws.Range("C8").NumberFormat = "0.000000"
$ws.Range("C8").WrapText = $true
$ws.Range("D8").Value = 0.00438273660367
$ws.Range("D8").NumberFormat = "0.000000"

# Row 9
$ws.Range("A9").Value = 4.4
$ws.Range("B9").Value = "RF 100, null=-1, balancing = 10"
$ws.Range("C9").Value = 0.843638300354
$ws.Range("C9").NumberFormat = "0.000000"
$ws.Range("C9").WrapText = $true
$ws.Range("D9").Value = 0.00264344442622
$ws.Range("D9").NumberFormat = "0.000000"

# Row heights
$ws.Rows("6").RowHeight = 14.9
$ws.Rows("7").RowHeight = 14.95
$ws.Rows("8").RowHeight = 14.95
$ws.Rows("9").RowHeight = 14.95
$ws.Rows("10").RowHeight = 13.8

# Selection matches the target workbook state
$ws.Range("E6").Select()
